$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Delete()
$paras = $d.Content.Paragraphs
$p5 = $paras.Item(5)
$r5 = $d.Range($p5.Range.Start, $p5.Range.End)
$r5.Delete()
$paras2 = $d.Content.Paragraphs
$p4 = $paras2.Item(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End)
$r4.Delete()

$paras3 = $d.Content.Paragraphs
$lastPara = $paras3.Item($paras3.Count)
$targetStart = $lastPara.Range.Start
Write-Host "target start" $targetStart

# Insert placeholder text
$insPoint = $d.Range($targetStart, $targetStart)
$insPoint.Text = "X"

Write-Host "Content End after insert:" $d.Content.End

# range covering the inserted X
$xRange = $d.Range($targetStart, $targetStart + 1)
Write-Host "xRange text=[$($xRange.Text)]"
$d.Bookmarks.Add("_GoBack", $xRange)

$bm = $d.Bookmarks.Item("_GoBack")
Write-Host "bm range" $bm.Range.Start $bm.Range.End

# Now delete the placeholder text
$xRange2 = $d.Range($targetStart, $targetStart + 1)
$xRange2.Text = ""
Write-Host "Content End after removing placeholder:" $d.Content.End

$bm2 = $d.Bookmarks.Item("_GoBack")
Write-Host "bm2 range" $bm2.Range.Start $bm2.Range.End
